# corrected data cleaning for pre/post/total fixation data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove header styling (bold font, thin border, centered alignment) ---
$ws.Range("A1:P1").ClearFormats()

# --- Clear the "Unnamed: 0" label in A1 ---
$ws.Range("A1").Value = ""

# --- Correct Revisit count row (row 3) ---
$ws.Range("B3").Value = 46
$ws.Range("D3").Value = 0
$ws.Range("H3").Value = 8
$ws.Range("K3").Value = 27
$ws.Range("L3").Value = 44

# --- Correct Fixation count row (row 4) ---
$ws.Range("B4").Value = 164
$ws.Range("D4").Value = 0
$ws.Range("H4").Value = 14
$ws.Range("K4").Value = 55
$ws.Range("L4").Value = 477

# --- Remove the "TTFF AOI (ms)" row entirely (row 5); remaining rows shift up ---
$ws.Rows("5:5").Delete()

# --- Remove the two trailing blank rows (now rows 9 and 10 after the shift) ---
$ws.Rows("9:10").Delete()

# --- Update values for "Dwell time (ms)" row (now row 5) ---
$ws.Range("B5").Value = 46242.84
$ws.Range("H5").Value = 6578.87
$ws.Range("K5").Value = 18501.87
$ws.Range("L5").Value = 153525.2

# --- Update values for "Dwell time (%)" row (now row 6) ---
$ws.Range("B6").Value = 14.09
$ws.Range("C6").Value = 1.2
$ws.Range("E6").Value = 0.42
$ws.Range("F6").Value = 0.97
$ws.Range("G6").Value = 1.11
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 0.62
$ws.Range("J6").Value = 0.23
$ws.Range("K6").Value = 5.64
$ws.Range("L6").Value = 46.78
$ws.Range("M6").Value = 0.31
$ws.Range("N6").Value = 0.05
$ws.Range("O6").Value = 0.72
$ws.Range("P6").Value = 0.24

# --- Update values for "Fixation duration (ms)" row (now row 7) ---
$ws.Range("B7").Value = 281.97
$ws.Range("H7").Value = 469.92
$ws.Range("K7").Value = 336.4
$ws.Range("L7").Value = 321.86

# Row 8 "First fixation duration (ms)" values are unchanged from the source row.

Write-Host ("Done. UsedRange: " + $ws.UsedRange.Address())
